# Reflection on the Honours Project - apply edits per commit "Completed draft of reflection"
#
# Strategy: use Range.InsertXML with a minimal pkg:package wrapper to splice exact
# OOXML run/paragraph fragments into the document. This gives precise control over
# run boundaries and special elements (w:lastRenderedPageBreak, bookmarks) that the
# higher level object model does not expose directly.

function Insert-Ooxml($range, [string]$fragment) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           $fragment +
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

function Replace-ParagraphContent($doc, [int]$index, [string]$fragment) {
    $p = $doc.Paragraphs.Item($index)
    $full = $p.Range
    $rng = $doc.Range($full.Start, $full.End - 1)
    Insert-Ooxml $rng $fragment
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# We work from the BOTTOM of the document upwards so that paragraph indices
# for content we have not processed yet remain valid (inserting/removing
# paragraphs below a given index never shifts that index).
# ---------------------------------------------------------------------------

# --- Change 8: replace the red "Overall critical evaluation..." paragraph (28)
#     and the trailing empty paragraph (29) with two new reflection paragraphs;
#     the _GoBack bookmark moves from paragraph 26 to the very end of the new
#     last paragraph.
$p28 = $d.Paragraphs.Item(28)
$p29 = $d.Paragraphs.Item(29)
$rng = $d.Range($p28.Range.Start, $p29.Range.End - 1)
$frag8 = '<w:body>' +
  '<w:p><w:r><w:t>While the completion of the project posed many technical challenges, the group dynamics were such that I was able to focus on the work and not on solving group issues. This allowed the final deliverable to be of a high quality and something that I am very proud of.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> The biggest lesson that I feel I can take away from this experience is that t</w:t></w:r>' +
  '<w:r><w:t>he relationship between the members of the team have a drastic impact on the final product.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> This is because when I compare this to my capstone project</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>of</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> 2013, </w:t></w:r>' +
  '<w:r><w:t>which was not completed to the expected quality, the only differe</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">nce is the relationship between group </w:t></w:r>' +
  '<w:r><w:t>members.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> While in both projects all team members where highly skilled,</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> the</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> honours project was completed to a higher quality because there was effective communication as w</w:t></w:r>' +
  '<w:r><w:t>ell as mutual trust and respect among group members.</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:r><w:t>Overall the honours projects and honours year has been a very rewarding experience and the lessons learned will continue to impact my future work.</w:t></w:r></w:p>' +
  '</w:body>'
Insert-Ooxml $rng $frag8

# --- Change 7: remove the _GoBack bookmark from its old location (end of the
#     "There are currently no plans..." paragraph, index 26).
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# Add the _GoBack bookmark back at the very end of the document (end of the
# new last paragraph, collapsed range right before its paragraph mark).
$last = $d.Paragraphs.Last
$endRng = $d.Range($last.Range.End - 1, $last.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $endRng) | Out-Null

# --- Change 6: tidy up run splits in the "made me feel incompetent" paragraph (21)
#     and drop the mid-paragraph lastRenderedPageBreak there.
$frag21 = '<w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">The fact that I had missed an important part of the system, made me feel incompetent for a few days </w:t></w:r>' +
  '<w:r><w:t>and I spent approximately a week trying alternative ways to fix it. After this I found a solution but realised that time would not allow me to implement it and complete the project. From this experience I learned an important lesson and it was</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">, to accept that I had made an error and </w:t></w:r>' +
  '<w:r><w:t>to try my best to fix it but if it is not possible then I should accept the consequences of my mistake and focus on completing the project to the best of my ability.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> It also taught me to consult the client at every step of the design process to verify that the design is </w:t></w:r>' +
  '<w:r><w:t>in line with their expectations</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> which again is something which I was aware of but did not have first-hand experience with.</w:t></w:r>' +
  '</w:p></w:body>'
Replace-ParagraphContent $d 21 $frag21

# --- Change 5: split the run in the "missed an important part" paragraph (20)
#     and add a lastRenderedPageBreak before "specific aspects of the project...".
$frag20 = '<w:body><w:p>' +
  '<w:r><w:t>During my design of the memo processing component of the project, I missed an important part that</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> I only discovered too late into discussion with the supervisor. The fact that I had missed this and that neither I nor Zahraa had noticed, made us realise that we had some communication issues which needed to be resolved. We realised that while we were communicating often, it was not very focused on </w:t></w:r>' +
  '<w:r><w:lastRenderedPageBreak/><w:t>specific aspects of the project and thus we made sure to resolve all concerns as soon as possible. This ensured that we would not forget, about it in future discussion.</w:t></w:r>' +
  '</w:p></w:body>'
Replace-ParagraphContent $d 20 $frag20

# --- Change 4: merge the trailing two runs in "...I feel the" / "project can be
#     deemed a success." (paragraph 15) into a single run.
$frag15 = '<w:body><w:p>' +
  '<w:r><w:t>The initial goals of the project were to create a web and tablet interfac</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">e for marking scripts and a web </w:t></w:r>' +
  '<w:r><w:t>application to view the test results. In the end a much more comprehensive test management solution was developed which limited the changes to the current system while still improving it in a meaningful way. Since more was accomplished than we initially set out to do and that all of this is of a high degree of quality, I feel the project can be deemed a success.</w:t></w:r>' +
  '</w:p></w:body>'
Replace-ParagraphContent $d 15 $frag15

# --- Change 3: remove the lastRenderedPageBreak from "Even with the motivation..." (14)
$p14 = $d.Paragraphs.Item(14)
$first14 = $d.Range($p14.Range.Start, $p14.Range.Start)
$frag14 = '<w:body><w:p><w:r><w:t>Even with the motivation of the user testing</w:t></w:r></w:p></w:body>'
# Replace only the first run (covering up to the end of "user testing") which
# currently carries the lastRenderedPageBreak.
$runEnd = $p14.Range.Start + [string]"Even with the motivation of the user testing".Length
$rngBreak = $d.Range($p14.Range.Start, $runEnd)
Insert-Ooxml $rngBreak '<w:body><w:p><w:r><w:t>Even with the motivation of the user testing</w:t></w:r></w:p></w:body>'

# --- Change 2: add a lastRenderedPageBreak before "In order to improve the..." (13)
$p13 = $d.Paragraphs.Item(13)
$runStart13 = $p13.Range.Start
$runEnd13 = $runStart13 + [string]"In order to improve the".Length
$rng13 = $d.Range($runStart13, $runEnd13)
Insert-Ooxml $rng13 '<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>In order to improve the</w:t></w:r></w:p></w:body>'

# --- Change 1: replace the empty paragraph after the "Introduction" heading (5)
#     with two new introductory paragraphs.
$frag5 = '<w:body>' +
  '<w:p>' +
  '<w:r><w:t xml:space="preserve">The purpose of this </w:t></w:r>' +
  '<w:r><w:t>report</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> is to r</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">eflect on the honours project. </w:t></w:r>' +
  '<w:r><w:t>To complete the project I worked with Zahraa Mathews. We knew each other prior to the project since we have been friends for three years and have worked in larger groups together as well.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> The project we completed, ScriptView, was chosen by us and</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> it</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> was our first choice.</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
  '<w:r><w:t xml:space="preserve">This report presents my experience of completing of the honours project </w:t></w:r>' +
  '<w:r><w:t>divided into the categories of computer science, project management, people and project impact.</w:t></w:r>' +
  '</w:p>' +
  '</w:body>'
Replace-ParagraphContent $d 5 $frag5

Write-Output "done"
